$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$win = $excel.ActiveWindow
$win.ScrollRow = 6
$win.ScrollColumn = 5
